# Update "想去人数" (column F) values on the 展览 and 全部类型 sheets
# to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 11611
    5  = 1055
    7  = 75
    9  = 48
    10 = 10964
    11 = 4222
    13 = 16
    15 = 2480
    16 = 1058
    17 = 65
    18 = 8
    19 = 144
    20 = 462
    21 = 11172
    22 = 10989
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
